# Daily attendance processing - 2025-10-08 11:41:01
# Normalize/reorder the "Recorded By" (column G) values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "System, backup@backdoor.com, system"
    3   = "System, dnasr281@gmail.com"
    4   = "System, backup@backdoor.com"
    6   = "System, dnasr281@gmail.com"
    11  = "System, dnasr281@gmail.com"
    12  = "System, dnasr281@gmail.com"
    13  = "System, dnasr281@gmail.com"
    14  = "System, dnasr281@gmail.com"
    15  = "System, dnasr281@gmail.com"
    29  = "System, backup@backdoor.com, system"
    30  = "System, dnasr281@gmail.com"
    33  = "System, dnasr281@gmail.com"
    38  = "System, dnasr281@gmail.com"
    39  = "System, dnasr281@gmail.com"
    40  = "System, dnasr281@gmail.com"
    41  = "System, dnasr281@gmail.com"
    42  = "System, dnasr281@gmail.com"
    56  = "System, backup@backdoor.com, system"
    57  = "System, dnasr281@gmail.com"
    58  = "System, backup@backdoor.com"
    60  = "System, dnasr281@gmail.com"
    65  = "System, dnasr281@gmail.com"
    66  = "System, dnasr281@gmail.com"
    67  = "System, dnasr281@gmail.com"
    68  = "System, dnasr281@gmail.com"
    69  = "System, dnasr281@gmail.com"
    86  = "System, dnasr281@gmail.com"
    89  = "System, dnasr281@gmail.com"
    90  = "dnasr281@gmail.com, admin@admin.com"
    93  = "System, dnasr281@gmail.com"
    95  = "System, dnasr281@gmail.com"
    112 = "System, dnasr281@gmail.com"
    115 = "System, dnasr281@gmail.com"
    116 = "dnasr281@gmail.com, admin@admin.com"
    119 = "System, dnasr281@gmail.com"
    121 = "System, dnasr281@gmail.com"
    138 = "System, dnasr281@gmail.com"
    141 = "System, dnasr281@gmail.com"
    142 = "dnasr281@gmail.com, admin@admin.com"
    145 = "System, dnasr281@gmail.com"
    147 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
